$d = $word.ActiveDocument

# Locate the last paragraph in the document ("install ads") and append
# two new list items after it, matching the existing list formatting.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara1.Range.Text = "input and .val"

$r2 = $newPara1.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara2.Range.Text = "where to put code for dropdown menu"
